# update with 16 apr cdc date corrected error in gr calculations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the "cases" (column B) values for 14-Mar through 14-Apr (rows 7-95) ---
$ws.Range("B7").Value2 = 1
$ws.Range("B8").Value2 = 2
$ws.Range("B9").Value2 = 2
$ws.Range("B10").Value2 = 5
$ws.Range("B11").Value2 = 2
$ws.Range("B18").Value2 = 3
$ws.Range("B20").Value2 = 7
$ws.Range("B22").Value2 = 10
$ws.Range("B23").Value2 = 9
$ws.Range("B24").Value2 = 4
$ws.Range("B26").Value2 = 4
$ws.Range("B29").Value2 = 5
$ws.Range("B30").Value2 = 5
$ws.Range("B33").Value2 = 11
$ws.Range("B35").Value2 = 17
$ws.Range("B36").Value2 = 16
$ws.Range("B37").Value2 = 12
$ws.Range("B39").Value2 = 26
$ws.Range("B40").Value2 = 33
$ws.Range("B41").Value2 = 41
$ws.Range("B42").Value2 = 47
$ws.Range("B43").Value2 = 36
$ws.Range("B44").Value2 = 62
$ws.Range("B45").Value2 = 109
$ws.Range("B46").Value2 = 92
$ws.Range("B47").Value2 = 126
$ws.Range("B48").Value2 = 128
$ws.Range("B49").Value2 = 207
$ws.Range("B50").Value2 = 170
$ws.Range("B51").Value2 = 441
$ws.Range("B52").Value2 = 413
$ws.Range("B53").Value2 = 466
$ws.Range("B54").Value2 = 530
$ws.Range("B55").Value2 = 600
$ws.Range("B56").Value2 = 794
$ws.Range("B57").Value2 = 946
$ws.Range("B58").Value2 = 1470
$ws.Range("B59").Value2 = 2505
$ws.Range("B60").Value2 = 2739
$ws.Range("B61").Value2 = 3543
$ws.Range("B62").Value2 = 5480
$ws.Range("B63").Value2 = 6916
$ws.Range("B64").Value2 = 7818
$ws.Range("B65").Value2 = 9496
$ws.Range("B66").Value2 = 11381
$ws.Range("B67").Value2 = 9497
$ws.Range("B68").Value2 = 9431
$ws.Range("B69").Value2 = 11087
$ws.Range("B70").Value2 = 13231
$ws.Range("B71").Value2 = 13020
$ws.Range("B72").Value2 = 13477
$ws.Range("B73").Value2 = 15208
$ws.Range("B74").Value2 = 11117
$ws.Range("B75").Value2 = 11505
$ws.Range("B76").Value2 = 16287
$ws.Range("B77").Value2 = 15455
$ws.Range("B78").Value2 = 14855
$ws.Range("B79").Value2 = 14753
$ws.Range("B80").Value2 = 15942
$ws.Range("B81").Value2 = 11648
$ws.Range("B82").Value2 = 11219
$ws.Range("B83").Value2 = 16092
$ws.Range("B84").Value2 = 15618
$ws.Range("B85").Value2 = 13968
$ws.Range("B86").Value2 = 11663
$ws.Range("B87").Value2 = 10886
$ws.Range("B88").Value2 = 6183
$ws.Range("B89").Value2 = 4176
$ws.Range("B90").Value2 = 4177
$ws.Range("B91").Value2 = 1363
$ws.Range("B92").Value2 = 339
$ws.Range("B93").Value2 = 182
$ws.Range("B94").Value2 = 102
$ws.Range("B95").Value2 = 19

# Row 85 (11-Apr) is no longer one of the most-recent/highlighted rows now that
# 16-Apr data has been added, so its yellow highlight is removed (copy the
# already-unhighlighted format from the row above, keeping A85's own value).
$ws.Range("A84").Copy()
$ws.Range("A85").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Append the new day: 16-Apr-2020 (serial 43936) with 0 cases so far ---
$ws.Range("A96").Value2 = 43936
$ws.Range("B96").Value2 = 0
$a96 = $ws.Range("A96")
$a96.NumberFormat = "[$-409]dd\-mmm\-yy;@"
$a96.Interior.Color = 65535

# --- Update the visible selection/scroll position to reflect where the user was working ---
$ws.Range("B75").Select() | Out-Null
